$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Header labels (rotate Black -> White -> Pink, Blue stays) ---
$ws.Range("B2").Value = "Uren Mr. White"
$ws.Range("E2").Value = "Uren Mr. Pink"
# H2 ("Uren Mr. Blue") is unchanged

# --- Legend labels in L column ---
$ws.Range("L5").Value = "Mr. White"
$ws.Range("L6").Value = "Mr. Pink"
# L7 ("Mr. Blue") is unchanged

# --- Data changes ---
$ws.Range("C4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("E5").Value = 41057

# --- Apply the date number format (matching B4/E4 style) down the B/E/H columns ---
# so the newly touched rows pick up the same "s=1" formatted-but-empty style
$ws.Range("B4").Copy()
$ws.Range("B4:B38").PasteSpecial(-4122)

$ws.Range("E4").Copy()
$ws.Range("E4:E36").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("H4:H37").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("F5").Select() | Out-Null
